# Re-order the SVI factor-variable lists (PCA variable grouping / labels)
# and refresh the dependent loading & variance statistics to match the
# re-ordered / re-computed factor analysis output.
$wb = $excel.ActiveWorkbook

# --- Sheet: Significant Components ---
$ws = $wb.Worksheets.Item('Significant Components')
$ws.Range('C2').Value = '[''QSERV'' ''QHISPC'' ''QEDLESHI'' ''QNOHLTH'' ''QESL'' ''PPUNIT'' ''QEXTRCT'' ''QFHH''
 ''PERCAP'']'
$ws.Range('C3').Value = '[''PERCAP'' ''QRICH'' ''MDHSEVAL'']'
$ws.Range('C4').Value = '[''QAGEDEP'' ''MEDAGE'' ''QSSBEN'']'
$ws.Range('C5').Value = '[''QAGEDEP'' ''QFEMLBR'' ''QFEMALE'']'

# --- Sheet: Loading Factors ---
$ws = $wb.Worksheets.Item('Loading Factors')
$ws.Range('A2').Value = 'QSERV'
$ws.Range('B2').Value = 0.5817049338713295
$ws.Range('C2').Value = 0.357728739689089
$ws.Range('D2').Value = -0.2244058845813555
$ws.Range('E2').Value = -0.03201466741222177
$ws.Range('F2').Value = 0.2754259411710515
$ws.Range('A3').Value = 'QHISPC'
$ws.Range('B3').Value = 0.8328587010890846
$ws.Range('C3').Value = 0.3339037451707143
$ws.Range('D3').Value = -0.1364957926512206
$ws.Range('E3').Value = -0.1269704285522055
$ws.Range('F3').Value = 0.09670565101208484
$ws.Range('A4').Value = 'QEDLESHI'
$ws.Range('B4').Value = 0.8777939327518004
$ws.Range('C4').Value = 0.2130779242463672
$ws.Range('D4').Value = -0.01839361158788122
$ws.Range('E4').Value = -0.1076883429830311
$ws.Range('F4').Value = 0.1846338453811155
$ws.Range('A5').Value = 'QNOHLTH'
$ws.Range('B5').Value = 0.6889886593829873
$ws.Range('C5').Value = 0.4149592263924404
$ws.Range('D5').Value = -0.1190897026747462
$ws.Range('E5').Value = -0.1154839644404644
$ws.Range('F5').Value = 0.2786568941960132
$ws.Range('A6').Value = 'QESL'
$ws.Range('B6').Value = 0.8009661451308956
$ws.Range('C6').Value = 0.1517950782154231
$ws.Range('D6').Value = -0.03424774610890587
$ws.Range('E6').Value = -0.2374898262054474
$ws.Range('F6').Value = 0.2037776753746855
$ws.Range('A7').Value = 'PPUNIT'
$ws.Range('B7').Value = 0.7302277632447152
$ws.Range('C7').Value = -0.004620267764485815
$ws.Range('D7').Value = -0.1512975013944436
$ws.Range('E7').Value = 0.05458729963722662
$ws.Range('F7').Value = -0.4752714245440762
$ws.Range('A8').Value = 'QEXTRCT'
$ws.Range('B8').Value = 0.7677511994868141
$ws.Range('C8').Value = 0.1449653467374824
$ws.Range('D8').Value = 0.01129678864827885
$ws.Range('E8').Value = -0.2382017668993546
$ws.Range('F8').Value = 0.09278788674142219
$ws.Range('B9').Value = 0.5630560198586153
$ws.Range('C9').Value = 0.3008187141192311
$ws.Range('D9').Value = -0.09551317420652185
$ws.Range('E9').Value = 0.2634165345840366
$ws.Range('F9').Value = -0.03178269667396563
$ws.Range('B10').Value = 0.4895374268553458
$ws.Range('C10').Value = 0.721454085175546
$ws.Range('D10').Value = -0.2685725771886172
$ws.Range('E10').Value = 0.05482089910664199
$ws.Range('F10').Value = 0.1831253495400197
$ws.Range('A11').Value = 'QRICH'
$ws.Range('B11').Value = 0.215097143564559
$ws.Range('C11').Value = 0.8701305635565615
$ws.Range('D11').Value = -0.1729034328402849
$ws.Range('E11').Value = -0.01425153590998154
$ws.Range('F11').Value = 0.2948729848375892
$ws.Range('A12').Value = 'MDHSEVAL'
$ws.Range('B12').Value = 0.3857664786261823
$ws.Range('C12').Value = 0.8013369886790841
$ws.Range('D12').Value = -0.03601197779313163
$ws.Range('E12').Value = -0.02873185420482495
$ws.Range('F12').Value = -0.03016150458661936
$ws.Range('A13').Value = 'QAGEDEP'
$ws.Range('B13').Value = -0.04282432248387342
$ws.Range('C13').Value = -0.1184872953833233
$ws.Range('D13').Value = 0.6543381063288032
$ws.Range('E13').Value = 0.6427560393897815
$ws.Range('F13').Value = -0.1139548187338235
$ws.Range('A14').Value = 'MEDAGE'
$ws.Range('B14').Value = -0.3105472636462374
$ws.Range('C14').Value = -0.2465384279923527
$ws.Range('D14').Value = 0.7910361825210296
$ws.Range('E14').Value = -0.01289474646594264
$ws.Range('F14').Value = -0.271708580213839
$ws.Range('B15').Value = 0.01836497920271595
$ws.Range('C15').Value = -0.05367045514883373
$ws.Range('D15').Value = 0.7773306974995353
$ws.Range('E15').Value = 0.1362099676011545
$ws.Range('F15').Value = -0.1455945088309715
$ws.Range('A16').Value = 'QFEMLBR'
$ws.Range('B16').Value = -0.2416846134945231
$ws.Range('C16').Value = 0.08178103491007524
$ws.Range('D16').Value = -0.02959553020525794
$ws.Range('E16').Value = 0.7849929108617704
$ws.Range('F16').Value = 0.003645845585976967
$ws.Range('A17').Value = 'QFEMALE'
$ws.Range('B17').Value = -0.04556235807336762
$ws.Range('C17').Value = -0.05659266919884228
$ws.Range('D17').Value = 0.1671554368380671
$ws.Range('E17').Value = 0.8778432099141861
$ws.Range('F17').Value = -0.02420501412752381
$ws.Range('B18').Value = 0.01759023568044128
$ws.Range('C18').Value = 0.2288188123825083
$ws.Range('D18').Value = -0.4234144335062181
$ws.Range('E18').Value = -0.09662324182207796
$ws.Range('F18').Value = 0.7659446312787375
$ws.Range('B19').Value = 0.1660837611468413
$ws.Range('C19').Value = 0.06295541642402984
$ws.Range('D19').Value = -0.106418739822054
$ws.Range('E19').Value = -0.0153032610956329
$ws.Range('F19').Value = 0.6312035431698786
$ws.Range('B20').Value = 0.370133440062485
$ws.Range('C20').Value = 0.1576334771461969
$ws.Range('D20').Value = -0.3817953717735523
$ws.Range('E20').Value = 0.08010852997758064
$ws.Range('F20').Value = 0.4611334617895158

# --- Sheet: All Refactor Variances ---
$ws = $wb.Worksheets.Item('All Refactor Variances')
$ws.Range('I2').Value = 4.839746419273361
$ws.Range('J2').Value = 3.422021257792846
$ws.Range('K2').Value = 2.232417330885144
$ws.Range('L2').Value = 2.05633789583397
$ws.Range('M2').Value = 2.044708482541559
$ws.Range('N2').Value = 5.117234768818557
$ws.Range('O2').Value = 2.680150598601196
$ws.Range('P2').Value = 2.243476584242993
$ws.Range('Q2').Value = 2.066288198356862
$ws.Range('R2').Value = 1.901714790296077
$ws.Range('I3').Value = 0.2304641152034934
$ws.Range('J3').Value = 0.1629533932282307
$ws.Range('K3').Value = 0.1063055871850068
$ws.Range('L3').Value = 0.09792085218257
$ws.Range('M3').Value = 0.0973670705972171
$ws.Range('N3').Value = 0.2693281457272925
$ws.Range('O3').Value = 0.1410605578211156
$ws.Range('P3').Value = 0.1180777149601576
$ws.Range('Q3').Value = 0.1087520104398348
$ws.Range('R3').Value = 0.1000902521208462
$ws.Range('I4').Value = 0.2304641152034934
$ws.Range('J4').Value = 0.3934175084317241
$ws.Range('K4').Value = 0.4997230956167309
$ws.Range('L4').Value = 0.5976439477993009
$ws.Range('M4').Value = 0.6950110183965179
$ws.Range('N4').Value = 0.2693281457272925
$ws.Range('O4').Value = 0.4103887035484081
$ws.Range('P4').Value = 0.5284664185085656
$ws.Range('Q4').Value = 0.6372184289484004
$ws.Range('R4').Value = 0.7373086810692466
$ws.Range('I5').Value = 0.3315977863706456
$ws.Range('J5').Value = 0.2344615968883281
$ws.Range('K5').Value = 0.1529552544796597
$ws.Range('L5').Value = 0.1408910788328023
$ws.Range('M5').Value = 0.1400942834285646
$ws.Range('N5').Value = 0.3652854667826672
$ws.Range('O5').Value = 0.1913181838799853
$ws.Range('P5').Value = 0.1601469208105905
$ws.Range('Q5').Value = 0.1474986165660256
$ws.Range('R5').Value = 0.1357508119607314

# --- Sheet: Final Variances ---
$ws = $wb.Worksheets.Item('Final Variances')
$ws.Range('B2').Value = 5.117234768818557
$ws.Range('C2').Value = 2.680150598601196
$ws.Range('D2').Value = 2.243476584242993
$ws.Range('E2').Value = 2.066288198356862
$ws.Range('F2').Value = 1.901714790296077
$ws.Range('B3').Value = 0.2693281457272925
$ws.Range('C3').Value = 0.1410605578211156
$ws.Range('D3').Value = 0.1180777149601576
$ws.Range('E3').Value = 0.1087520104398348
$ws.Range('F3').Value = 0.1000902521208462
$ws.Range('B4').Value = 0.2693281457272925
$ws.Range('C4').Value = 0.4103887035484081
$ws.Range('D4').Value = 0.5284664185085656
$ws.Range('E4').Value = 0.6372184289484004
$ws.Range('F4').Value = 0.7373086810692466
$ws.Range('B5').Value = 0.3652854667826672
$ws.Range('C5').Value = 0.1913181838799853
$ws.Range('D5').Value = 0.1601469208105905
$ws.Range('E5').Value = 0.1474986165660256
$ws.Range('F5').Value = 0.1357508119607314

# --- Sheet: Included and Excluded ---
$ws = $wb.Worksheets.Item('Included and Excluded')
$ws.Range('B2').Value = '[[''QSERV'', ''QHISPC'', ''QEDLESHI'', ''QNOHLTH'', ''QESL'', ''PPUNIT'', ''QEXTRCT'', ''QFHH'', ''PERCAP'', ''QRICH'', ''MDHSEVAL'', ''QAGEDEP'', ''MEDAGE'', ''QSSBEN'', ''QFEMLBR'', ''QFEMALE'', ''QRENTER'', ''QNOAUTO'', ''QPOVTY'']]'
